$wb = $excel.ActiveWorkbook

# ALC!row12
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 627.8
$ws.Range("I12").Value = 586.44446
$ws.Range("K12").Value = 586.44446
$ws.Range("M12").Value = -416.44446

# ALC!row41
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 599.8
$ws.Range("I41").Value = 445
$ws.Range("J41").Value = 703
$ws.Range("K41").Value = 445
$ws.Range("L41").Value = 703
$ws.Range("M41").Value = -5
$ws.Range("N41").Value = -1583

# ALC!row104
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H104").Value = 230
$ws.Range("I104").Value = 285
$ws.Range("J104").Value = 120
$ws.Range("K104").Value = 855
$ws.Range("L104").Value = 360
$ws.Range("M104").Value = 892
$ws.Range("N104").Value = -3854

# ALC!row112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1999.75
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").ClearContents()

# ALC!row137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2551.3333
$ws.Range("I137").Value = 1500
$ws.Range("K137").Value = 4500
$ws.Range("M137").Value = -1950

# ARM!row41
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 1991.3334
$ws.Range("I41").Value = 693.3333
$ws.Range("J41").Value = 4587.3335
$ws.Range("K41").Value = 693.3333
$ws.Range("L41").Value = 4587.3335
$ws.Range("M41").Value = -279.3333
$ws.Range("N41").Value = -5415.3335

# ARM!row44
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

# ARM!row97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 30305438
$ws.Range("I97").Value = 33335482
$ws.Range("K97").Value = 33335482
$ws.Range("M97").Value = -33334986

# ARM!row135
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 41283.2
$ws.Range("J135").Value = 41283.2
$ws.Range("L135").Value = 41283.2
$ws.Range("N135").Value = -51423.2

# BSM!row20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1639.75
$ws.Range("I20").Value = 1525
$ws.Range("J20").Value = 1754.5
$ws.Range("K20").Value = 1525
$ws.Range("L20").Value = 1754.5
$ws.Range("M20").Value = -1278
$ws.Range("N20").Value = -2248.5

# BSM!row86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3129.4443
$ws.Range("I86").Value = 3129.4443
$ws.Range("K86").Value = 3129.4443
$ws.Range("M86").Value = -2006.4443

# BSM!row89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3129.4443
$ws.Range("I89").Value = 3129.4443
$ws.Range("K89").Value = 15647.2215
$ws.Range("M89").Value = -10031.2215

# BSM!row94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 159993.28
$ws.Range("I94").Value = 368903
$ws.Range("J94").Value = 3311
$ws.Range("K94").Value = 368903
$ws.Range("L94").Value = 3311
$ws.Range("M94").Value = -368452
$ws.Range("N94").Value = -4213

# BSM!row134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2993.5715
$ws.Range("I134").Value = 2993
$ws.Range("K134").Value = 8979
$ws.Range("M134").Value = -6444

# CRP!row47
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H47").Value = 19750
$ws.Range("I47").Value = 2500
$ws.Range("K47").Value = 2500
$ws.Range("M47").Value = -1934

# CRP!row141
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 855184.25
$ws.Range("J141").Value = 855184.25
$ws.Range("L141").Value = 855184.25
$ws.Range("N141").Value = -865544.25

# CUL!row14
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 2582.4546
$ws.Range("I14").Value = 2582.4546
$ws.Range("K14").Value = 7747.3638
$ws.Range("M14").Value = -7574.3638

# CUL!row80
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2250.5
$ws.Range("J80").Value = 2302
$ws.Range("L80").Value = 6906
$ws.Range("N80").Value = -8778

# CUL!row83
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 2250.5
$ws.Range("J83").Value = 2302
$ws.Range("L83").Value = 20718
$ws.Range("N83").Value = -30078

# CUL!row107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 700.3333
$ws.Range("I107").Value = 280
$ws.Range("J107").Value = 1036.6
$ws.Range("K107").Value = 840
$ws.Range("L107").Value = 3109.8
$ws.Range("M107").Value = 1080
$ws.Range("N107").Value = -6949.799999999999

# GSM!row22
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 4186.875
$ws.Range("J22").Value = 4186.875
$ws.Range("L22").Value = 4186.875
$ws.Range("N22").Value = -5244.875

# GSM!row58
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 51347
$ws.Range("I58").Value = 50041
$ws.Range("J58").Value = 52000
$ws.Range("K58").Value = 50041
$ws.Range("L58").Value = 52000
$ws.Range("M58").Value = -49764
$ws.Range("N58").Value = -52554

# GSM!row80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3600
$ws.Range("I80").Value = 2750
$ws.Range("J80").Value = 4166.6665
$ws.Range("K80").Value = 2750
$ws.Range("L80").Value = 4166.6665
$ws.Range("M80").Value = -1752
$ws.Range("N80").Value = -6162.6665

# GSM!row83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3600
$ws.Range("I83").Value = 2750
$ws.Range("J83").Value = 4166.6665
$ws.Range("K83").Value = 13750
$ws.Range("L83").Value = 20833.3325
$ws.Range("M83").Value = -8758
$ws.Range("N83").Value = -30817.3325

# GSM!row102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2291
$ws.Range("I102").Value = 1612
$ws.Range("J102").Value = 2970
$ws.Range("K102").Value = 1612
$ws.Range("L102").Value = 2970
$ws.Range("M102").Value = 10
$ws.Range("N102").Value = -6214

# GSM!row122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2670.077
$ws.Range("I122").Value = 2554.111
$ws.Range("K122").Value = 7662.333
$ws.Range("M122").Value = -5212.333

# LTW!row25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 852
$ws.Range("I25").Value = 800
$ws.Range("J25").Value = 904
$ws.Range("K25").Value = 800
$ws.Range("L25").Value = 904
$ws.Range("M25").Value = -570
$ws.Range("N25").Value = -1364

# LTW!row55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1078
$ws.Range("I55").Value = 840.2
$ws.Range("J55").Value = 1315.8
$ws.Range("K55").Value = 840.2
$ws.Range("L55").Value = 1315.8
$ws.Range("M55").Value = -667.2
$ws.Range("N55").Value = -1661.8

# LTW!row57
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

# LTW!row127
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 10000
$ws.Range("J127").Value = 10000
$ws.Range("L127").Value = 10000
$ws.Range("N127").Value = -19920

# WVR!row32
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 26068.467
$ws.Range("I32").Value = 18205.4
$ws.Range("K32").Value = 18205.4
$ws.Range("M32").Value = -17888.4

# WVR!row43
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 35000
$ws.Range("J43").Value = 35000
$ws.Range("L43").Value = 35000
$ws.Range("N43").Value = -35298

# WVR!row126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2735.4285
$ws.Range("I126").Value = 1529.6
$ws.Range("K126").Value = 4588.799999999999
$ws.Range("M126").Value = -2118.799999999999

# WVR!row136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3358.1428
$ws.Range("I136").Value = 3141.4
$ws.Range("K136").Value = 9424.200000000001
$ws.Range("M136").Value = -6874.200000000001
